# Commit: "Added some sample sketches"
# - Renames the "Pole Number" register to "Pole ID" and updates its factory
#   default value from 0x00 to 0xFF.
# - Adds a new "Time" register row (Regular Registers section) right after
#   the existing Config Registers block.
# - Adds a new "<! Regular Registers>" section header a few rows below.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update existing "Pole" register row (row 23) ---
$ws.Range("A23").Value = "Pole ID"
$ws.Range("F23").Value = "0xFF (Factory Default)"

# --- New register row describing a 4-byte unix time value ---
$ws.Range("A25").Value = "Time"
$ws.Range("B25").Value = 13
$ws.Range("C25").Value = 4
$ws.Range("D25").Value = "R/W"
$ws.Range("E25").Value = "4 Byte unix time(Integer)"

# --- New section header further down the sheet ---
$ws.Range("A28").Value = "<! Regular Registers>"

# Match the author's final selection/active cell
[void]$ws.Range("E23").Select()
